$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "have a disadvantage economically" -> "have an advantage economically"
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "have a disadvantage economically",
    $false, $true, $false, $false, $false, $true, 1, $false,
    "have an advantage economically", 2)

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from its old spot (in the middle of
#    "overblo|wn, and most programs...") down to the end of the CONCLUSION
#    paragraph, right after "...rarely fixes the underlying cause."
#    (Bookmarks.Add with the existing name re-seats the single bookmark,
#    which both removes it from the old location and creates it at the
#    new one.)
# ---------------------------------------------------------------------------
$text = $d.Content.Text
$marker = "rarely fixes the underlying cause."
$idx = $text.IndexOf($marker)
$endPos = $idx + $marker.Length
$bmRange = $d.Range($endPos, $endPos)
$null = $d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 3) Rename section header "WHAT IS NEEDED INSTEAD" ->
#    "INTERNET ACCESS IS READILY AVAILABLE"
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "WHAT IS NEEDED INSTEAD",
    $false, $true, $false, $false, $false, $true, 1, $false,
    "INTERNET ACCESS IS READILY AVAILABLE", 2)

# ---------------------------------------------------------------------------
# 4) In the CONCLUSION paragraph, the word break around the page-break run
#    shifts from "access" / " do not" to "access do" / " not" (the
#    surrounding gramStart/gramEnd proofing marks around "access" are also
#    gone in the final copy). Re-split the run boundary accordingly; the
#    lastRenderedPageBreak stays attached to the run that now ends in
#    "access do".
# ---------------------------------------------------------------------------
$text = $d.Content.Text
$needle = "Programs that provide free Internet access do not address this issue."
$idx = $text.IndexOf($needle)
if ($idx -ge 0) {
    $splitPos = $idx + "Programs that provide free Internet access do".Length
    $r = $d.Range($splitPos, $splitPos)
    $r.InsertAfter("|")
}

Write-Output "done"
